$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: Manufacturer part number + price corrections
$ws.Range("C12").Value = "CL21B104KCFNNNE"
$ws.Range("G12").Value = 0.015
$ws.Range("H12").Value = 0.015

# Row 13: Manufacturer part number changes to numeric value, price corrections
$ws.Range("C13").Value = 885012207127
$ws.Range("G13").Value = 0.1
$ws.Range("H13").Value = 0.1

# Row 25: Footprint correction
$ws.Range("E25").Value = "6-0805_M"

# Row 38: Manufacturer part number + price corrections
$ws.Range("C38").Value = "RC0805FR-07100RL"
$ws.Range("G38").Value = 0.006
$ws.Range("H38").Value = 0.006
